$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-id"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")
# Keep the "Fixed Value" for Extension.url in sync with the new URL
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-id"
# The top-level Extension row no longer carries the ele-1/ext-1 constraint text
$elem.Range("AI2").Value = ""
